# "Update pawn support panel"
#
# - Delete event rows 401-405 (rows 24-28), the "Dangerous_Mission_2_x" /
#   chapter-4 combat-log block that was dropped from the design doc.
# - Normalize the whitespace in the three "combat tutorial" description
#   strings referenced by rows 10-12 (C10:C12): the double blank lines
#   collapse to a single line break, and the stray leading space before
#   "Different roles..." is removed.
# - Rows 10-12 therefore wrap onto two lines instead of three/four, so
#   their row height shrinks to 27.6pt.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the trailing "Dangerous_Mission" rows (401-405) entirely - this
# shifts everything below back up and shrinks the sheet dimension from
# A1:F28 to A1:F23.
$ws.Range("A24:F28").EntireRow.Delete()

# Re-write the three tutorial blurbs with normalized whitespace.
$ws.Range("C10").Value = "When your units enterd a area with controled by an enemy,  they will start a combat. `nUnits will constantly consume ammo duing combat.Their damage will decrease once they ran out of ammo."
$ws.Range("C11").Value = "You can spend bullets/medecine to set an ammo box/med kit in a area, which will benefit your units in combat. Click on any areas controled by your unit, and then press the construction button, choose a item you want to set."
$ws.Range("C12").Value = "When a unit is in a combat, units that stand in the areas next to it will offer support to the combat unit (as the green line shows). Different roles will have different support skills."

# The shorter text now only wraps across two lines, so these rows get
# shorter.
$ws.Rows.Item(10).RowHeight = 27.6
$ws.Rows.Item(11).RowHeight = 27.6
$ws.Rows.Item(12).RowHeight = 27.6
